# Generate Report for Handback
# Updates the handoff/handback timestamps for the
# "f3e2fe82-073d-4902-8608-044de7c6793f.md" file row on each sheet, as a
# new localization round trip was generated for it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the f3e2fe82 row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-02 18:54:22"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H3) / "Correspond Handback DateTime" (K3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-02 18:54:15"
$wsZhCn.Range("K3").Value = "2016-09-02 18:54:37"

# --- de-de sheet: "Correspond Handoff Datetime" (H3) / "Correspond Handback DateTime" (K3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-02 18:54:22"
$wsDeDe.Range("K3").Value = "2016-09-02 18:54:45"
